# word_freq_list.xlsx — retire four stimulus words and add their
# replacements (plus refreshed frequency counts), matching the commit's
# "Added new trial lists" cleanup of the natural/artificial word lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A12: כוורת -> שבתאי ; frequency 21 -> 18 -----------------------------
$ws.Range("A12").Value = "שבתאי"
$ws.Range("B12").Value = 18

# --- A25: צללים -> ירקות ; frequency 12 -> 77 -----------------------------
# (this row previously carried the yellow "flag" highlight; the edit
# clears it along with the word swap)
$ws.Range("A25").Clear()
$ws.Range("A25").Value = "ירקות"
$ws.Range("A25").Borders.Item(7).Color = 0
$ws.Range("A25").Borders.Item(7).LineStyle = 1
$ws.Range("B25").Value = 77

# --- C28: גלידה -> צוללת ; frequency 27 -> 25 -----------------------------
$ws.Range("C28").Value = "צוללת"
$ws.Range("D28").Value = 25

# --- A21: פיסגה -> תמרים ; frequency 18 -> 16 -----------------------------
# (also had the yellow highlight, cleared the same way as A25)
$ws.Range("A21").Clear()
$ws.Range("A21").Value = "תמרים"
$ws.Range("A21").Borders.Item(7).Color = 0
$ws.Range("A21").Borders.Item(7).LineStyle = 1
$ws.Range("B21").Value = 16

# --- restore the selection left behind in the saved workbook -------------
$ws.Range("B22").Select()

$wb.Save()
